# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E64) is reordered from reverse-chronological
# (2007 down to 1607) to chronological order (1607 up to 2007), and the
# "Valor Mora" column (F16:F64) is realigned so that the value actually
# tracks the period it sits next to (27578 for 1607-1808, 31249 for
# 1809-2006, 22916 for 2007 only).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $period = $periods[$i]

    $ws.Cells.Item($row, 5).Value = $period

    if ($period -eq "2007") {
        $value = 22916
    } elseif ([int]$period -le 1808) {
        $value = 27578
    } else {
        $value = 31249
    }

    $ws.Cells.Item($row, 6).Value = $value
}
